$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vm_pu result values: slack/reference voltage setpoint (column B)
# changed from 1.05 to 1.02 p.u., with recalculated voltage magnitudes
# across the other buses (columns C-F, I-N) for the 380 kV case.

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.069229783573056
$ws.Cells.Item(2, 4).Value = 1.073452509572986
$ws.Cells.Item(2, 5).Value = 1.082210962593434
$ws.Cells.Item(2, 6).Value = 1.087900031501516
$ws.Cells.Item(2, 9).Value = 1.054513651681027
$ws.Cells.Item(2, 10).Value = 1.074165567252045
$ws.Cells.Item(2, 11).Value = 1.076144332288658
$ws.Cells.Item(2, 12).Value = 1.084879771175308
$ws.Cells.Item(2, 13).Value = 1.090554105452439
$ws.Cells.Item(2, 14).Value = 1.028268844460397

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.070391972735182
$ws.Cells.Item(3, 4).Value = 1.074403302909823
$ws.Cells.Item(3, 5).Value = 1.083334106665573
$ws.Cells.Item(3, 6).Value = 1.089035133506138
$ws.Cells.Item(3, 9).Value = 1.054864250002646
$ws.Cells.Item(3, 10).Value = 1.074983822671005
$ws.Cells.Item(3, 11).Value = 1.076911733837148
$ws.Cells.Item(3, 12).Value = 1.08582076818021
$ws.Cells.Item(3, 13).Value = 1.091508098975943
$ws.Cells.Item(3, 14).Value = 1.028551202615609

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.071144077049202
$ws.Cells.Item(4, 4).Value = 1.075018571266204
$ws.Cells.Item(4, 5).Value = 1.084061274338114
$ws.Cells.Item(4, 6).Value = 1.089770071559804
$ws.Cells.Item(4, 9).Value = 1.055089950305015
$ws.Cells.Item(4, 10).Value = 1.075512792434533
$ws.Cells.Item(4, 11).Value = 1.077407708074269
$ws.Cells.Item(4, 12).Value = 1.086429479284547
$ws.Cells.Item(4, 13).Value = 1.092125253488765
$ws.Cells.Item(4, 14).Value = 1.028733533588351

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.071460284241375
$ws.Cells.Item(5, 4).Value = 1.075277240183781
$ws.Cells.Item(5, 5).Value = 1.084367076207081
$ws.Cells.Item(5, 6).Value = 1.090079147923228
$ws.Cells.Item(5, 9).Value = 1.055184557097266
$ws.Cells.Item(5, 10).Value = 1.0757350529211
$ws.Cells.Item(5, 11).Value = 1.077616075669626
$ws.Cells.Item(5, 12).Value = 1.086685339489801
$ws.Cells.Item(5, 13).Value = 1.092384671365606
$ws.Cells.Item(5, 14).Value = 1.028810096072822

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.071513378169474
$ws.Cells.Item(6, 4).Value = 1.075320672430464
$ws.Cells.Item(6, 5).Value = 1.084418427596908
$ws.Cells.Item(6, 6).Value = 1.090131049573262
$ws.Cells.Item(6, 9).Value = 1.055200425718993
$ws.Cells.Item(6, 10).Value = 1.075772364503068
$ws.Cells.Item(6, 11).Value = 1.077651053306001
$ws.Cells.Item(6, 12).Value = 1.086728297096312
$ws.Cells.Item(6, 13).Value = 1.092428226792839
$ws.Cells.Item(6, 14).Value = 1.02882294600487

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.071148302135423
$ws.Cells.Item(7, 4).Value = 1.075022027573035
$ws.Cells.Item(7, 5).Value = 1.084065360081903
$ws.Cells.Item(7, 6).Value = 1.089774201026829
$ws.Cells.Item(7, 9).Value = 1.055091215536242
$ws.Cells.Item(7, 10).Value = 1.075515762754774
$ws.Cells.Item(7, 11).Value = 1.077410492841912
$ws.Cells.Item(7, 12).Value = 1.086432898265427
$ws.Cells.Item(7, 13).Value = 1.092128719975779
$ws.Cells.Item(7, 14).Value = 1.028734556971526

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.069622532728148
$ws.Cells.Item(8, 4).Value = 1.073773826143135
$ws.Cells.Item(8, 5).Value = 1.082590448022931
$ws.Cells.Item(8, 6).Value = 1.088283551363382
$ws.Cells.Item(8, 9).Value = 1.054632378355372
$ws.Cells.Item(8, 10).Value = 1.074442203480069
$ws.Cells.Item(8, 11).Value = 1.076403800717515
$ws.Cells.Item(8, 12).Value = 1.085197822385043
$ws.Cells.Item(8, 13).Value = 1.090876541883316
$ws.Cells.Item(8, 14).Value = 1.028364346080022

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.06693458067245
$ws.Cells.Item(9, 4).Value = 1.071574645521977
$ws.Cells.Item(9, 5).Value = 1.079994649240101
$ws.Cells.Item(9, 6).Value = 1.085660273262468
$ws.Cells.Item(9, 9).Value = 1.053814955447302
$ws.Cells.Item(9, 10).Value = 1.072546636691924
$ws.Cells.Item(9, 11).Value = 1.074625380257627
$ws.Cells.Item(9, 12).Value = 1.083020087816883
$ws.Cells.Item(9, 13).Value = 1.088668931039165
$ws.Cells.Item(9, 14).Value = 1.027709125724138

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.06514298703034
$ws.Cells.Item(10, 4).Value = 1.07010871941003
$ws.Cells.Item(10, 5).Value = 1.078266227960562
$ws.Cells.Item(10, 6).Value = 1.083913703659388
$ws.Cells.Item(10, 9).Value = 1.053264016473825
$ws.Cells.Item(10, 10).Value = 1.071280331144992
$ws.Cells.Item(10, 11).Value = 1.073436720128075
$ws.Cells.Item(10, 12).Value = 1.08156731336182
$ws.Cells.Item(10, 13).Value = 1.087196415501956
$ws.Cells.Item(10, 14).Value = 1.027270387363691

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.064367281491755
$ws.Cells.Item(11, 4).Value = 1.06947399930937
$ws.Cells.Item(11, 5).Value = 1.077518295650315
$ws.Cells.Item(11, 6).Value = 1.083157953866322
$ws.Cells.Item(11, 9).Value = 1.053024030078302
$ws.Cells.Item(11, 10).Value = 1.070731384722
$ws.Cells.Item(11, 11).Value = 1.072921289232543
$ws.Cells.Item(11, 12).Value = 1.080938013228559
$ws.Cells.Item(11, 13).Value = 1.086558609453185
$ws.Cells.Item(11, 14).Value = 1.027079951729694

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.06407915853064
$ws.Cells.Item(12, 4).Value = 1.069238240900603
$ws.Cells.Item(12, 5).Value = 1.077240552473874
$ws.Cells.Item(12, 6).Value = 1.082877313094762
$ws.Cells.Item(12, 9).Value = 1.052934673862356
$ws.Cells.Item(12, 10).Value = 1.070527386583868
$ws.Cells.Item(12, 11).Value = 1.072729724535431
$ws.Cells.Item(12, 12).Value = 1.080704226490004
$ws.Cells.Item(12, 13).Value = 1.086321669592536
$ws.Cells.Item(12, 14).Value = 1.027009146401827

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.064140961515721
$ws.Cells.Item(13, 4).Value = 1.069288811665158
$ws.Cells.Item(13, 5).Value = 1.077300126073659
$ws.Cells.Item(13, 6).Value = 1.082937507959072
$ws.Cells.Item(13, 9).Value = 1.052953850789391
$ws.Cells.Item(13, 10).Value = 1.070571149186613
$ws.Cells.Item(13, 11).Value = 1.072770820832152
$ws.Cells.Item(13, 12).Value = 1.080754376208911
$ws.Cells.Item(13, 13).Value = 1.086372495378869
$ws.Cells.Item(13, 14).Value = 1.027024337511317

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.064343464989988
$ws.Cells.Item(14, 4).Value = 1.069454511348439
$ws.Cells.Item(14, 5).Value = 1.077495335844526
$ws.Cells.Item(14, 6).Value = 1.083134754418473
$ws.Cells.Item(14, 9).Value = 1.053016648244867
$ws.Cells.Item(14, 10).Value = 1.070714524114059
$ws.Cells.Item(14, 11).Value = 1.072905456696642
$ws.Cells.Item(14, 12).Value = 1.080918689095568
$ws.Cells.Item(14, 13).Value = 1.086539024552364
$ws.Cells.Item(14, 14).Value = 1.027074100348415

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.064468235165037
$ws.Cells.Item(15, 4).Value = 1.069556605027664
$ws.Cells.Item(15, 5).Value = 1.077615620569715
$ws.Cells.Item(15, 6).Value = 1.083256294832564
$ws.Cells.Item(15, 9).Value = 1.053055311386006
$ws.Cells.Item(15, 10).Value = 1.070802849530775
$ws.Cells.Item(15, 11).Value = 1.072988395609159
$ws.Cells.Item(15, 12).Value = 1.081019922810384
$ws.Cells.Item(15, 13).Value = 1.086641624637407
$ws.Cells.Item(15, 14).Value = 1.027104751721665

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.065194469334333
$ws.Cells.Item(16, 4).Value = 1.070150844404103
$ws.Cells.Item(16, 5).Value = 1.078315875929231
$ws.Cells.Item(16, 6).Value = 1.08396387131186
$ws.Cells.Item(16, 9).Value = 1.05327991350551
$ws.Cells.Item(16, 10).Value = 1.071316749626802
$ws.Cells.Item(16, 11).Value = 1.073470912083626
$ws.Cells.Item(16, 12).Value = 1.081609072838556
$ws.Cells.Item(16, 13).Value = 1.087238740377759
$ws.Cells.Item(16, 14).Value = 1.027283016270623

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.065650033646336
$ws.Cells.Item(17, 4).Value = 1.070523604275136
$ws.Cells.Item(17, 5).Value = 1.078755257066142
$ws.Cells.Item(17, 6).Value = 1.084407855752632
$ws.Cells.Item(17, 9).Value = 1.05342041842821
$ws.Cells.Item(17, 10).Value = 1.071638937029248
$ws.Cells.Item(17, 11).Value = 1.073773385225288
$ws.Cells.Item(17, 12).Value = 1.081978566541368
$ws.Cells.Item(17, 13).Value = 1.087613242053484
$ws.Cells.Item(17, 14).Value = 1.027394713953904

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.065915763092526
$ws.Cells.Item(18, 4).Value = 1.070741032264298
$ws.Cells.Item(18, 5).Value = 1.079011587544981
$ws.Cells.Item(18, 6).Value = 1.084666875201409
$ws.Cells.Item(18, 9).Value = 1.053502234987407
$ws.Cells.Item(18, 10).Value = 1.071826802965498
$ws.Cells.Item(18, 11).Value = 1.073949741983957
$ws.Cells.Item(18, 12).Value = 1.082194062934518
$ws.Cells.Item(18, 13).Value = 1.087831663572986
$ws.Cells.Item(18, 14).Value = 1.027459821039273

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.066006371108652
$ws.Cells.Item(19, 4).Value = 1.070815170190755
$ws.Cells.Item(19, 5).Value = 1.079098997591398
$ws.Cells.Item(19, 6).Value = 1.084755202784666
$ws.Cells.Item(19, 9).Value = 1.053530108992562
$ws.Cells.Item(19, 10).Value = 1.071890850149547
$ws.Cells.Item(19, 11).Value = 1.074009863104901
$ws.Cells.Item(19, 12).Value = 1.082267537750613
$ws.Cells.Item(19, 13).Value = 1.087906136469447
$ws.Cells.Item(19, 14).Value = 1.027482013364354

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.065601155245287
$ws.Cells.Item(20, 4).Value = 1.07048361031364
$ws.Cells.Item(20, 5).Value = 1.0787081107831
$ws.Cells.Item(20, 6).Value = 1.084360215170853
$ws.Cells.Item(20, 9).Value = 1.053405357813002
$ws.Cells.Item(20, 10).Value = 1.071604375635233
$ws.Cells.Item(20, 11).Value = 1.073740940037248
$ws.Cells.Item(20, 12).Value = 1.081938925751239
$ws.Cells.Item(20, 13).Value = 1.08757306353506
$ws.Cells.Item(20, 14).Value = 1.027382734430642

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.064283832537129
$ws.Cells.Item(21, 4).Value = 1.069405716790071
$ws.Cells.Item(21, 5).Value = 1.07743784944589
$ws.Cells.Item(21, 6).Value = 1.083076668095396
$ws.Cells.Item(21, 9).Value = 1.052998161885485
$ws.Cells.Item(21, 10).Value = 1.070672306399231
$ws.Cells.Item(21, 11).Value = 1.072865812850937
$ws.Cells.Item(21, 12).Value = 1.080870304063251
$ws.Cells.Item(21, 13).Value = 1.086489986703195
$ws.Cells.Item(21, 14).Value = 1.027059448340159

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.063455628638114
$ws.Cells.Item(22, 4).Value = 1.068728030464124
$ws.Cells.Item(22, 5).Value = 1.076639602994155
$ws.Cells.Item(22, 6).Value = 1.08227010406115
$ws.Cells.Item(22, 9).Value = 1.052740899977054
$ws.Cells.Item(22, 10).Value = 1.070085727530385
$ws.Cells.Item(22, 11).Value = 1.072314944829871
$ws.Cells.Item(22, 12).Value = 1.080198206800657
$ws.Cells.Item(22, 13).Value = 1.085808837692054
$ws.Cells.Item(22, 14).Value = 1.026855785972227

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.063894670984894
$ws.Cells.Item(23, 4).Value = 1.069087282292577
$ws.Cells.Item(23, 5).Value = 1.077062729242114
$ws.Cells.Item(23, 6).Value = 1.082697636225056
$ws.Cells.Item(23, 9).Value = 1.052877397182798
$ws.Cells.Item(23, 10).Value = 1.070396736442697
$ws.Cells.Item(23, 11).Value = 1.072607031318809
$ws.Cells.Item(23, 12).Value = 1.080554518645878
$ws.Cells.Item(23, 13).Value = 1.086169944545315
$ws.Cells.Item(23, 14).Value = 1.026963789151981

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.065623241270323
$ws.Cells.Item(24, 4).Value = 1.070501681852408
$ws.Cells.Item(24, 5).Value = 1.078729414013839
$ws.Cells.Item(24, 6).Value = 1.084381741741743
$ws.Cells.Item(24, 9).Value = 1.053412163481724
$ws.Cells.Item(24, 10).Value = 1.071619992628274
$ws.Cells.Item(24, 11).Value = 1.07375560083939
$ws.Cells.Item(24, 12).Value = 1.081956837788884
$ws.Cells.Item(24, 13).Value = 1.087591218536433
$ws.Cells.Item(24, 14).Value = 1.027388147598317

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.067629410208684
$ws.Cells.Item(25, 4).Value = 1.072143151153373
$ws.Cells.Item(25, 5).Value = 1.080665351130733
$ws.Cells.Item(25, 6).Value = 1.086338049150253
$ws.Cells.Item(25, 9).Value = 1.054027333897763
$ws.Cells.Item(25, 10).Value = 1.07303714142387
$ws.Cells.Item(25, 11).Value = 1.075085680022428
$ws.Cells.Item(25, 12).Value = 1.083583250028595
$ws.Cells.Item(25, 13).Value = 1.089239786004963
$ws.Cells.Item(25, 14).Value = 1.027878855080541
